$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work from the bottom rows upward so we don't clobber values we still need to read/copy.

# Row 7 (new): step 6 text + "denied access" result (F7 keeps its existing style/blank cell)
$ws.Range("C7").Value = "Step 6: Try to delete an assessment of someone who is not under my district"
$ws.Range("D7").Value = "I am denied access to this"

# Row 6 (new): step 5 text (renumbered from old step 4) + "denied access" result (unchanged text)
$ws.Range("C6").Value = "Step 5: While logged in try to delete an assessment about me"
$ws.Range("D6").Value = "I am denied access to this"

# Row 5 (new): step 4 text (renumbered from old step 3) + "data removed" result (unchanged text)
$ws.Range("C5").Value = "Step 4: Delete one that belongs to someone from my district"
$ws.Range("D5").Value = "The data is removed from the database."

# Row 4 (new): brand new step 3 content
$ws.Range("C4").Value = 'Step 3: Go to the "Assessments" page'
$ws.Range("D4").Value = "I am redirected to the assessments page"

# Row 3 (new): rewritten step 2 content
$ws.Range("C3").Value = "Step 2: Log in as a user with the appropriate role"
$ws.Range("D3").Value = "I am redirected to the user's dashboard"

# Update the active selection to D4 to match the saved view state
$ws.Range("D4").Select()
